# Refresh the cryptos price/volume snapshot (GitHub Actions scheduled update).
# Price ("D") and Volume(1h) ("E") columns hold plain text (inline strings) in
# the source workbook, e.g. "37.487.06" (dotted thousands, no real decimal
# meaning) or "  +0.87%  " (padded percent). Excel's COM layer auto-detects
# numeric-looking text on assignment, so for column D we force the cell to
# Text format before writing the new value and then clear the format again
# (NumberFormat "@" -> write -> ClearFormats) so the value is stored as a
# string while the cell keeps its original (default) style, matching the
# target file exactly. Column E values always contain a "%" and spaces, so
# they are never mis-detected as numbers and can be assigned directly.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '37.494.73'
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = '  +0.84%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.016.37'
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = '  +0.77%  '

$ws.Range("E4").Value = '  -0.02%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '263.52'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  +7.03%  '

$ws.Range("E6").Value = '  -1.61%  '

$ws.Range("E7").Value = '  +0.00%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '56.08'
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  -6.93%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.385'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  +0.22%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0772'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  -3.88%  '

$ws.Range("E11").Value = '  -2.17%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '14.37'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  -4.26%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '2.313.36'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  +0.90%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.806'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  -4.80%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '20.87'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  -8.09%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '5.25'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  -3.90%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.035.12'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  +1.60%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '37.378.19'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  +0.85%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '69.70'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  -0.84%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.0₃0842'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  -2.45%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '5.17'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  -0.27%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '228.03'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  -1.11%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '2.68'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  +9.17%  '

$ws.Range("E24").Value = '  -0.04%  '

$ws.Range("E25").Value = '  -0.60%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '165.03'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  +0.89%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.01'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  -4.65%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '19.68'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  +0.02%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.128'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  -10.44%  '

$ws.Range("E30").Value = '  -2.33%  '

$ws.Range("E31").Value = '  -1.14%  '

$ws.Range("B32").Value = 'Filecoin'
$ws.Range("C32").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.64'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  -3.47%  '

$ws.Range("B33").Value = 'Hedera'
$ws.Range("C33").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.0649'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  -1.62%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.53'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  -0.09%  '

$ws.Range("E35").Value = '  +0.61%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.83'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  +1.28%  '

$ws.Range("E37").Value = '  +0.11%  '

$ws.Range("E38").Value = '  +1.42%  '

$ws.Range("E39").Value = '  -4.50%  '

$ws.Range("E40").Value = '  +4.12%  '

$ws.Range("E41").Value = '  +3.43%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.0938'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  -4.30%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.0213'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  -0.84%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.393.03'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  +1.42%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '90.26'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  -1.04%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '15.69'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  -6.04%  '

$ws.Range("E47").Value = '  -1.73%  '

$ws.Range("B48").Value = 'FraxShare'
$ws.Range("C48").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '7.05'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  -4.09%  '

$ws.Range("B49").Value = 'MXToken'
$ws.Range("C49").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.91'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  +2.16%  '

$ws.Range("B50").Value = 'NEARProtocol'
$ws.Range("C50").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.97'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  -3.07%  '

$ws.Range("B51").Value = 'RocketPoolETH'
$ws.Range("C51").Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.205.62'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  +0.86%  '
